$d = $word.ActiveDocument

$replacements = @(
  @("719÷4=179, 3", "135÷7=19, 2"),
  @("920÷2=460, 0", "452÷3=150, 2"),
  @("112÷9=12, 4", "910÷3=303, 1"),
  @("562÷5=112, 2", "101÷2=50, 1"),
  @("929÷3=309, 2", "823÷3=274, 1"),
  @("468÷7=66, 6", "544÷6=90, 4"),
  @("746÷5=149, 1", "625÷9=69, 4"),
  @("279÷5=55, 4", "458÷5=91, 3"),
  @("522÷5=104, 2", "356÷4=89, 0"),
  @("483÷5=96, 3", "895÷5=179, 0"),
  @("355÷5=71, 0", "741÷3=247, 0"),
  @("546÷8=68, 2", "766÷9=85, 1"),
  @("654÷6=109, 0", "778÷8=97, 2"),
  @("462÷9=51, 3", "143÷9=15, 8"),
  @("935÷9=103, 8", "188÷2=94, 0"),
  @("303÷7=43, 2", "242÷2=121, 0"),
  @("342÷9=38, 0", "724÷9=80, 4"),
  @("135÷5=27, 0", "510÷5=102, 0"),
  @("586÷8=73, 2", "512÷6=85, 2"),
  @("659÷5=131, 4", "496÷4=124, 0"),
  @("205÷9=22, 7", "793÷3=264, 1"),
  @("942÷8=117, 6", "193÷8=24, 1"),
  @("569÷2=284, 1", "842÷7=120, 2"),
  @("830÷8=103, 6", "912÷8=114, 0"),
  @("804÷4=201, 0", "385÷2=192, 1")
)

foreach ($pair in $replacements) {
  $old = $pair[0]
  $new = $pair[1]
  $range = $d.Content
  $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
